$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B48:F48").Copy()
$ws.Range("B49:F49").PasteSpecial(-4122)

$ws.Range("B49").Value = 39
$ws.Range("E49").Value = "https://programmingport.hashnode.dev/if-statement-or-shell-scripting"
$ws.Range("C49").Value = "If Statement | Shell Scripting "
$ws.Range("D49").Value = 44168
$ws.Range("F49").Value = "https://dev.to/rahulmishra05/if-statement-shell-scripting-434j"

$ws.ListObjects.Item("Table2").Resize($ws.Range("B10:F49"))
